# The document has two consecutive "signature block" paragraphs that each
# contain only a manual line break (<w:br/>) right before the signer's name
# ("A T E N T A M E N T E" ... blank ... blank ... "CLAUDIA GISELA ...").
# The commit adds one extra, completely empty paragraph (no runs at all)
# between those two line-break paragraphs, re-using the same paragraph
# formatting (spacing after=0/line=240/auto, bold Garet paragraph mark),
# in order to give the template a bit more vertical space.

$d = $word.ActiveDocument

# Locate the paragraph that consists solely of a manual line break and is
# immediately followed by "A T E N T A M E N T E" (the first of the two
# consecutive break-only paragraphs).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text
    if ($text -eq [string][char]11 + [string][char]13) {
        $prevText = $d.Paragraphs($i - 1).Range.Text
        if ($prevText -like "*A T E N T A M E N T E*") {
            $targetIndex = $i
            break
        }
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the target paragraph"
}

$targetPara = $d.Paragraphs($targetIndex)
$paraEnd = $targetPara.Range.End

# Range covering just the paragraph mark at the end of the target paragraph.
$markRange = $d.Range($paraEnd - 1, $paraEnd)

# Insert a new paragraph mark right before the existing one. This splits
# off a brand-new, completely empty paragraph (no runs) that inherits the
# paragraph formatting/properties of the original paragraph, matching the
# diff exactly.
$markRange.InsertBefore([string][char]13)
